$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.49"
$ws.Range("E2").Value = "'-1.58%"
$ws.Range("D3").Value = "'35.74"
$ws.Range("E3").Value = "'-0.61%"
$ws.Range("D4").Value = "'5.033"
$ws.Range("E4").Value = "'-0.66%"
$ws.Range("D5").Value = "'0.07970"
$ws.Range("E5").Value = "'-2.00%"
$ws.Range("D6").Value = "'1.857"
$ws.Range("E6").Value = "'-4.63%"
$ws.Range("D7").Value = "'4.122"
$ws.Range("E7").Value = "'-0.41%"
$ws.Range("D8").Value = "'7.762"
$ws.Range("E8").Value = "'-0.50%"
$ws.Range("D9").Value = "'0.9229"
$ws.Range("E9").Value = "'-1.41%"
$ws.Range("D10").Value = "'0.1269"
$ws.Range("E10").Value = "'-4.15%"
$ws.Range("D11").Value = "'0.1888"
$ws.Range("E11").Value = "'-1.63%"
$ws.Range("D12").Value = "'0.08963"
$ws.Range("E12").Value = "'-3.27%"
$ws.Range("D13").Value = "'0.03417"
$ws.Range("E13").Value = "'-2.62%"
$ws.Range("E14").Value = "'-0.27%"
$ws.Range("D15").Value = "'0.001400"
$ws.Range("E15").Value = "'-3.08%"
$ws.Range("D16").Value = "'0.006271"
$ws.Range("E16").Value = "'9.06%"
$ws.Range("D17").Value = "'3.861"
$ws.Range("E17").Value = "'7.09%"
$ws.Range("D18").Value = "'3.336"
$ws.Range("E18").Value = "'13.38%"
$ws.Range("D19").Value = "'0.3405"
$ws.Range("E19").Value = "'-0.68%"
$ws.Range("D20").Value = "'0.1340"
$ws.Range("E20").Value = "'0.62%"
$ws.Range("D21").Value = "'4.803"
$ws.Range("E21").Value = "'-7.18%"
$ws.Range("D22").Value = "'0.2342"
$ws.Range("E22").Value = "'-10.41%"
$ws.Range("D23").Value = "'0.04355"
$ws.Range("E23").Value = "'-0.70%"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'1.13%"
$ws.Range("D25").Value = "'0.004840"
$ws.Range("E25").Value = "'1.27%"
$ws.Range("E27").Value = "'-21.13%"
$ws.Range("E28").Value = "'42.31%"
$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "'-3.76%"
$ws.Range("D40").Value = "'0.05109"
$ws.Range("E40").Value = "'0.99%"
$ws.Range("D41").Value = "'0.007565"
$ws.Range("E41").Value = "'-0.75%"
$ws.Range("D42").Value = "'0.01012"
$ws.Range("E42").Value = "'-9.81%"
$ws.Range("D43").Value = "'0.1346"
$ws.Range("E43").Value = "'-2.58%"
$ws.Range("E44").Value = "'0.73%"
$ws.Range("D45").Value = "'0.009862"
$ws.Range("E45").Value = "'-12.68%"
$ws.Range("D46").Value = "'0.00006199"
$ws.Range("E46").Value = "'-2.94%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.24%"
$ws.Range("D48").Value = "'64.85"
$ws.Range("E48").Value = "'-0.17%"
$ws.Range("D49").Value = "'0.001252"
$ws.Range("E49").Value = "'5.25%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.24%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.24%"
